$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '22.258.92'
$ws.Range("E2").Value = '  -1.28%  '

# Row 3
$ws.Range("D3").Value = '1.558.07'
$ws.Range("E3").Value = '  -1.20%  '

# Row 4
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  -0.16%  '

# Row 5
$ws.Range("D5").Value = '''1.001'
$ws.Range("E5").Value = '  -0.07%  '

# Row 6
$ws.Range("D6").Value = '''288.73'
$ws.Range("E6").Value = '  -0.20%  '

# Row 7
$ws.Range("D7").Value = '''0.3799'
$ws.Range("E7").Value = '  +1.73%  '

# Row 8
$ws.Range("D8").Value = '''0.3298'
$ws.Range("E8").Value = '  -1.73%  '

# Row 9
$ws.Range("E9").Value = '  -7.90%  '

# Row 10
$ws.Range("D10").Value = '''1.139'
$ws.Range("E10").Value = '  +0.02%  '

# Row 11
$ws.Range("D11").Value = '''0.07385'
$ws.Range("E11").Value = '  -1.76%  '

# Row 12
$ws.Range("D12").Value = '''1.001'
$ws.Range("E12").Value = '  -0.14%  '

# Row 13
$ws.Range("D13").Value = '''20.24'
$ws.Range("E13").Value = '  -3.84%  '

# Row 14
$ws.Range("D14").Value = '''5.851'
$ws.Range("E14").Value = '  -2.45%  '

# Row 15
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '''6.759'
$ws.Range("E15").Value = '  -2.81%  '

# Row 16
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '1.548.78'
$ws.Range("E16").Value = '  -2.19%  '

# Row 17
$ws.Range("E17").Value = '  -4.24%  '

# Row 18
$ws.Range("D18").Value = '''0.06662'
$ws.Range("E18").Value = '  -1.69%  '

# Row 19
$ws.Range("D19").Value = '''86.54'
$ws.Range("E19").Value = '  -2.49%  '

# Row 20
$ws.Range("D20").Value = '''6.416'
$ws.Range("E20").Value = '  +0.06%  '

# Row 21
$ws.Range("D21").Value = '''1.002'
$ws.Range("E21").Value = '  +0.04%  '

# Row 22
$ws.Range("D22").Value = '''16.17'
$ws.Range("E22").Value = '  -2.33%  '

# Row 23
$ws.Range("D23").Value = '''11.73'
$ws.Range("E23").Value = '  -3.45%  '

# Row 24
$ws.Range("D24").Value = '22.254.12'
$ws.Range("E24").Value = '  -1.29%  '

# Row 25
$ws.Range("D25").Value = '''2.278'
$ws.Range("E25").Value = '  -5.22%  '

# Row 26
$ws.Range("E26").Value = '  -0.98%  '

# Row 27
$ws.Range("D27").Value = '''150.88'
$ws.Range("E27").Value = '  -1.45%  '

# Row 29
$ws.Range("D29").Value = '''4.937'
$ws.Range("E29").Value = '  -1.43%  '

# Row 30
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").Value = '''122.82'
$ws.Range("E30").Value = '  -1.33%  '

# Row 31
$ws.Range("B31").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C31").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D31").Value = '1.727.69'
$ws.Range("E31").Value = '  -1.80%  '

# Row 32
$ws.Range("D32").Value = '''1.091'
$ws.Range("E32").Value = '  +3.48%  '

# Row 33
$ws.Range("D33").Value = '''5.920'
$ws.Range("E33").Value = '  -4.57%  '

# Row 34
$ws.Range("D34").Value = '''1.922'
$ws.Range("E34").Value = '  -4.65%  '

# Row 35
$ws.Range("D35").Value = '''9.386'
$ws.Range("E35").Value = '  -3.89%  '

# Row 36
$ws.Range("D36").Value = '''0.08216'
$ws.Range("E36").Value = '  -1.48%  '

# Row 37
$ws.Range("D37").Value = '''0.02351'
$ws.Range("E37").Value = '  -4.65%  '

# Row 38
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").Value = '''5.351'
$ws.Range("E38").Value = '  -1.33%  '

# Row 39
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '''0.06316'
$ws.Range("E39").Value = '  -1.53%  '

# Row 40
$ws.Range("D40").Value = '''0.2156'
$ws.Range("E40").Value = '  -5.93%  '

# Row 41
$ws.Range("D41").Value = '''1.232'
$ws.Range("E41").Value = '  -5.28%  '

# Row 42
$ws.Range("E42").Value = '  -2.97%  '

# Row 43
$ws.Range("D43").Value = '''0.6076'
$ws.Range("E43").Value = '  -4.10%  '

# Row 44
$ws.Range("E44").Value = '  -0.05%  '

# Row 45
$ws.Range("D45").Value = '''13.73'
$ws.Range("E45").Value = '  -1.79%  '

# Row 46
$ws.Range("D46").Value = '''3.749'
$ws.Range("E46").Value = '  -1.09%  '

# Row 47
$ws.Range("D47").Value = '''0.5892'
$ws.Range("E47").Value = '  -4.34%  '

# Row 48
$ws.Range("D48").Value = '''123.03'
$ws.Range("E48").Value = '  -2.04%  '

# Row 49
$ws.Range("D49").Value = '''1.971'
$ws.Range("E49").Value = '  -4.80%  '

# Row 50
$ws.Range("E50").Value = '  -3.67%  '

# Row 51
$ws.Range("D51").Value = '''0.07065'
$ws.Range("E51").Value = '  -3.06%  '
